$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.368.56'
$ws.Range('E2').Value = '  -2.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.386.67'
$ws.Range('E3').Value = '  -2.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.37'
$ws.Range('E5').Value = '  -1.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.06'
$ws.Range('E6').Value = '  -4.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.387.19'
$ws.Range('E8').Value = '  -2.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.467'
$ws.Range('E9').Value = '  -3.42%  '
$ws.Range('E10').Value = '  +4.58%  '
$ws.Range('E11').Value = '  -6.63%  '
$ws.Range('E12').Value = '  -4.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.958.56'
$ws.Range('E13').Value = '  -2.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000199'
$ws.Range('E14').Value = '  -7.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.52'
$ws.Range('E15').Value = '  -7.17%  '
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.379.02'
$ws.Range('E17').Value = '  -2.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.384.09'
$ws.Range('E18').Value = '  -2.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.40'
$ws.Range('E19').Value = '  +2.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.10'
$ws.Range('E20').Value = '  -5.71%  '
$ws.Range('E21').Value = '  -5.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '413.93'
$ws.Range('E22').Value = '  -6.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.577'
$ws.Range('E23').Value = '  -5.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '76.85'
$ws.Range('E24').Value = '  -2.77%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.516.98'
$ws.Range('E26').Value = '  -2.53%  '
$ws.Range('E27').Value = '  -10.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.19'
$ws.Range('E28').Value = '  -6.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.73'
$ws.Range('E29').Value = '  -7.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.41'
$ws.Range('E30').Value = '  -3.05%  '
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.159'
$ws.Range('E32').Value = '  -4.66%  '
$ws.Range('E33').Value = '  -9.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '24.29'
$ws.Range('E34').Value = '  -4.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.382.50'
$ws.Range('E35').Value = '  -2.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.68'
$ws.Range('E37').Value = '  -7.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.51'
$ws.Range('E38').Value = '  -9.18%  '
$ws.Range('E39').Value = '  -5.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.21%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '167.49'
$ws.Range('E41').Value = '  -4.14%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0853'
$ws.Range('E42').Value = '  -4.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.868'
$ws.Range('E43').Value = '  -1.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.01'
$ws.Range('E44').Value = '  -7.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.90'
$ws.Range('E45').Value = '  -11.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '45.26'
$ws.Range('E46').Value = '  -1.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.60'
$ws.Range('E47').Value = '  -9.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.17'
$ws.Range('E48').Value = '  -5.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.01'
$ws.Range('E49').Value = '  -6.19%  '
$ws.Range('E50').Value = '  -9.18%  '
$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.232'
$ws.Range('E51').Value = '  -6.06%  '
